$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.815.71"
$ws.Range("E2").Value = "  +6.36%  "
$ws.Range("D3").Value = "2.304.02"
$ws.Range("E3").Value = "  +3.27%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.50%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "2.656.67"
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").Value = "2.303.84"
$ws.Range("E15").Value = "  +3.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.821"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.42%  "
$ws.Range("D18").Value = "46.831.17"
$ws.Range("E18").Value = "  +6.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +19.74%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "44.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0803"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.115"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.41%  "
$ws.Range("E37").Value = "  +3.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +22.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.06%  "
$ws.Range("E41").Value = "  +7.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0306"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.81%  "
$ws.Range("D45").Value = "1.858.18"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +20.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.198"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "75.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.39%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "96.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.34%  "

Write-Output "Applied 100 cell updates"
